$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "2026-02-06 Friday" "2026-02-07 Saturday"

Replace-Text "792÷8=99, 0" "749÷6=124, 5"
Replace-Text "223÷3=74, 1" "843÷5=168, 3"
Replace-Text "511÷9=56, 7" "463÷3=154, 1"
Replace-Text "567÷7=81, 0" "373÷6=62, 1"
Replace-Text "451÷2=225, 1" "667÷2=333, 1"

Replace-Text "534÷7=76, 2" "896÷2=448, 0"
Replace-Text "562÷8=70, 2" "305÷3=101, 2"
Replace-Text "347÷7=49, 4" "169÷5=33, 4"
Replace-Text "356÷9=39, 5" "753÷6=125, 3"
Replace-Text "354÷5=70, 4" "830÷3=276, 2"

Replace-Text "973÷7=139, 0" "419÷5=83, 4"
Replace-Text "317÷5=63, 2" "544÷7=77, 5"
Replace-Text "264÷2=132, 0" "457÷2=228, 1"
Replace-Text "756÷3=252, 0" "299÷7=42, 5"
Replace-Text "292÷7=41, 5" "544÷6=90, 4"

Replace-Text "663÷8=82, 7" "119÷9=13, 2"
Replace-Text "347÷4=86, 3" "988÷9=109, 7"
Replace-Text "584÷3=194, 2" "462÷5=92, 2"
Replace-Text "309÷4=77, 1" "183÷9=20, 3"
Replace-Text "663÷6=110, 3" "961÷7=137, 2"

Replace-Text "431÷4=107, 3" "494÷4=123, 2"
Replace-Text "260÷6=43, 2" "567÷2=283, 1"
Replace-Text "167÷3=55, 2" "783÷5=156, 3"
Replace-Text "833÷6=138, 5" "102÷8=12, 6"
Replace-Text "765÷6=127, 3" "230÷7=32, 6"
